# Update the cryptos price/volume table (columns D and E) for rows 2-51
# to reflect the latest scrape values from the GitHub Actions update.
#
# Price column (D) values that look purely numeric are written with a
# leading apostrophe so Excel stores them as literal text (preserving
# trailing zeros / exact formatting) instead of converting to a Double.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.503.80"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "1.677.88"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'219.95"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").Value = "'0.5307"
$ws.Range("E6").Value = "  +1.62%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.2696"
$ws.Range("E8").Value = "  +3.51%  "
$ws.Range("D9").Value = "'0.06416"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("D10").Value = "'21.79"
$ws.Range("E10").Value = "  +4.88%  "
$ws.Range("D11").Value = "'0.07802"
$ws.Range("E11").Value = "  +1.52%  "
$ws.Range("D12").Value = "1.683.96"
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").Value = "'4.511"
$ws.Range("E13").Value = "  +2.07%  "
$ws.Range("D14").Value = "'0.5597"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "0.0₅8365"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").Value = "'65.74"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "26.529.46"
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'4.797"
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("D20").Value = "'193.32"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").Value = "'10.34"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("D22").Value = "'6.328"
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'0.1272"
$ws.Range("E24").Value = "  +5.32%  "
$ws.Range("D25").Value = "'138.77"
$ws.Range("E25").Value = "  -4.89%  "
$ws.Range("D26").Value = "'7.418"
$ws.Range("D27").Value = "'16.31"
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("D28").Value = "'1.441"
$ws.Range("E28").Value = "  +3.48%  "
$ws.Range("D29").Value = "'0.06328"
$ws.Range("E29").Value = "  +6.89%  "
$ws.Range("D30").Value = "'1.290"
$ws.Range("E30").Value = "  +2.08%  "
$ws.Range("D31").Value = "'3.607"
$ws.Range("E31").Value = "  +4.96%  "
$ws.Range("D32").Value = "'3.444"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").Value = "'1.694"
$ws.Range("E33").Value = "  +2.47%  "
$ws.Range("D34").Value = "'1.013"
$ws.Range("E34").Value = "  +2.84%  "
$ws.Range("D35").Value = "'0.6180"
$ws.Range("E35").Value = "  +8.95%  "
$ws.Range("D36").Value = "'2.424"
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("D37").Value = "'2.784"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("D38").Value = "'0.01633"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("D39").Value = "'6.117"
$ws.Range("E39").Value = "  +5.83%  "
$ws.Range("D40").Value = "1.095.44"
$ws.Range("E40").Value = "  +6.58%  "
$ws.Range("D41").Value = "'0.8625"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "'100.70"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").Value = "1.823.39"
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("D45").Value = "0.0₈111"
$ws.Range("E45").Value = "  +2.47%  "
$ws.Range("D46").Value = "'58.73"
$ws.Range("E46").Value = "  +4.87%  "
$ws.Range("D47").Value = "'8.196"
$ws.Range("E47").Value = "  +1.55%  "
$ws.Range("D48").Value = "'0.9973"
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").Value = "'1.489"
$ws.Range("E49").Value = "  +7.46%  "
$ws.Range("D50").Value = "'0.05197"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").Value = "'6.035"
$ws.Range("E51").Value = "  +1.76%  "
